$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "empty fields" claim-creation test case to reflect that it is only
# an attempted creation (the text used to describe a different, now-removed
# scenario for the same row).
$ws.Range("A19").Value = "Попытка создания новой заявки с пустыми полями"

# Mark rows 18-22 (the first five "Заявка" test cases) as automated by copying
# the existing "V" formatting from column D of row 4 (which already carries the
# correct style) into the previously-empty D18:D22 cells.
for ($r = 18; $r -le 22; $r++) {
    $ws.Range("D4").Copy($ws.Range("D$r"))
}

# Restore the selected cell as saved in the workbook.
$ws.Range("D22").Select()
